$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1474.7273
$ws.Cells.Item(12, 10).Value = 1499.4
$ws.Cells.Item(12, 12).Value = 1499.4
$ws.Cells.Item(12, 14).Value = -1839.4
$ws.Cells.Item(15, 8).Value = 1223.5952
$ws.Cells.Item(15, 9).Value = 1223.5952
$ws.Cells.Item(15, 11).Value = 3670.7856
$ws.Cells.Item(15, 13).Value = -3501.7856
$ws.Cells.Item(18, 8).Value = 1739.3043
$ws.Cells.Item(18, 9).Value = 1718.3636
$ws.Cells.Item(18, 11).Value = 1718.3636
$ws.Cells.Item(18, 13).Value = -1434.3636
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 14).Value = 0
$ws.Cells.Item(29, 8).Value = 93
$ws.Cells.Item(29, 9).Value = 93
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 279
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 2
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(58, 8).Value = 6288.778
$ws.Cells.Item(58, 9).Value = 219.8
$ws.Cells.Item(58, 10).Value = 13875
$ws.Cells.Item(58, 11).Value = 659.4000000000001
$ws.Cells.Item(58, 12).Value = 41625
$ws.Cells.Item(58, 13).Value = -509.4000000000001
$ws.Cells.Item(58, 14).Value = -41925
$ws.Cells.Item(59, 8).Value = 2980
$ws.Cells.Item(59, 10).Value = 2980
$ws.Cells.Item(59, 12).Value = 8940
$ws.Cells.Item(59, 14).Value = -10054
$ws.Cells.Item(69, 8).Value = 16999.75
$ws.Cells.Item(69, 9).Value = 16999.75
$ws.Cells.Item(69, 11).Value = 50999.25
$ws.Cells.Item(69, 13).Value = -50125.25
$ws.Cells.Item(72, 8).Value = 16999.75
$ws.Cells.Item(72, 9).Value = 16999.75
$ws.Cells.Item(72, 11).Value = 152997.75
$ws.Cells.Item(72, 13).Value = -148629.75
$ws.Cells.Item(93, 8).Value = 49999
$ws.Cells.Item(93, 10).Value = 49999
$ws.Cells.Item(93, 12).Value = 49999
$ws.Cells.Item(93, 14).Value = -54991
$ws.Cells.Item(94, 8).Value = 9994.5
$ws.Cells.Item(94, 9).Value = 9994.5
$ws.Cells.Item(94, 11).Value = 9994.5
$ws.Cells.Item(94, 13).Value = -9543.5
$ws.Cells.Item(97, 8).Value = 2146.5
$ws.Cells.Item(97, 10).Value = 2410.2856
$ws.Cells.Item(97, 12).Value = 7230.8568
$ws.Cells.Item(97, 14).Value = -8222.856800000001
$ws.Cells.Item(113, 8).Value = 5029.25
$ws.Cells.Item(113, 9).Value = 4327.3076
$ws.Cells.Item(113, 10).Value = 5637.6
$ws.Cells.Item(113, 11).Value = 4327.3076
$ws.Cells.Item(113, 12).Value = 5637.6
$ws.Cells.Item(113, 13).Value = -1073.3076
$ws.Cells.Item(113, 14).Value = -12145.6
$ws.Cells.Item(121, 8).Value = 1746.7
$ws.Cells.Item(121, 10).Value = 1746.7
$ws.Cells.Item(121, 12).Value = 5240.1
$ws.Cells.Item(121, 14).Value = -8734.1
$ws.Cells.Item(132, 8).Value = 3788.2273
$ws.Cells.Item(132, 9).Value = 3450.4443
$ws.Cells.Item(132, 11).Value = 10351.3329
$ws.Cells.Item(132, 13).Value = -7821.332900000001
$ws.Cells.Item(138, 8).Value = 3984.2888
$ws.Cells.Item(138, 9).Value = 4638.8335
$ws.Cells.Item(138, 10).Value = 3746.2727
$ws.Cells.Item(138, 11).Value = 13916.5005
$ws.Cells.Item(138, 12).Value = 11238.8181
$ws.Cells.Item(138, 13).Value = -8776.500499999998
$ws.Cells.Item(138, 14).Value = -21518.8181
$ws.Cells.Item(141, 8).Value = 8809.111000000001
$ws.Cells.Item(141, 9).Value = 5345.3335
$ws.Cells.Item(141, 10).Value = 15736.667
$ws.Cells.Item(141, 11).Value = 16036.0005
$ws.Cells.Item(141, 12).Value = 47210.001
$ws.Cells.Item(141, 13).Value = -10856.0005
$ws.Cells.Item(141, 14).Value = -57570.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 970.75
$ws.Cells.Item(2, 9).Value = 952.875
$ws.Cells.Item(2, 10).Value = 1006.5
$ws.Cells.Item(2, 11).Value = 952.875
$ws.Cells.Item(2, 12).Value = 1006.5
$ws.Cells.Item(2, 13).Value = -839.875
$ws.Cells.Item(2, 14).Value = -1232.5
$ws.Cells.Item(32, 8).Value = 869.64935
$ws.Cells.Item(32, 9).Value = 674.01495
$ws.Cells.Item(32, 10).Value = 2180.4
$ws.Cells.Item(32, 11).Value = 674.01495
$ws.Cells.Item(32, 12).Value = 2180.4
$ws.Cells.Item(32, 13).Value = -387.01495
$ws.Cells.Item(32, 14).Value = -2754.4
$ws.Cells.Item(37, 8).Value = 54662.668
$ws.Cells.Item(37, 10).Value = 54994
$ws.Cells.Item(37, 12).Value = 54994
$ws.Cells.Item(37, 14).Value = -55540
$ws.Cells.Item(61, 8).Value = 2692.5715
$ws.Cells.Item(61, 9).Value = 1974.75
$ws.Cells.Item(61, 11).Value = 1974.75
$ws.Cells.Item(61, 13).Value = -1762.75
$ws.Cells.Item(63, 8).Value = 2559.111
$ws.Cells.Item(63, 10).Value = 2180.3333
$ws.Cells.Item(63, 12).Value = 2180.3333
$ws.Cells.Item(63, 14).Value = -3552.3333
$ws.Cells.Item(66, 8).Value = 2559.111
$ws.Cells.Item(66, 10).Value = 2180.3333
$ws.Cells.Item(66, 12).Value = 10901.6665
$ws.Cells.Item(66, 14).Value = -17765.6665
$ws.Cells.Item(74, 8).Value = 1329.5927
$ws.Cells.Item(74, 9).Value = 1174.421
$ws.Cells.Item(74, 10).Value = 1698.125
$ws.Cells.Item(74, 11).Value = 1174.421
$ws.Cells.Item(74, 12).Value = 1698.125
$ws.Cells.Item(74, 13).Value = -300.421
$ws.Cells.Item(74, 14).Value = -3446.125
$ws.Cells.Item(77, 8).Value = 1329.5927
$ws.Cells.Item(77, 9).Value = 1174.421
$ws.Cells.Item(77, 10).Value = 1698.125
$ws.Cells.Item(77, 11).Value = 5872.105
$ws.Cells.Item(77, 12).Value = 8490.625
$ws.Cells.Item(77, 13).Value = -1504.105
$ws.Cells.Item(77, 14).Value = -17226.625
$ws.Cells.Item(110, 8).Value = 4000
$ws.Cells.Item(110, 9).Value = 3000
$ws.Cells.Item(110, 10).Value = 5000
$ws.Cells.Item(110, 11).Value = 3000
$ws.Cells.Item(110, 12).Value = 5000
$ws.Cells.Item(110, 13).Value = -955
$ws.Cells.Item(110, 14).Value = -9090
$ws.Cells.Item(116, 8).Value = 970.75
$ws.Cells.Item(116, 9).Value = 952.875
$ws.Cells.Item(116, 10).Value = 1006.5
$ws.Cells.Item(116, 11).Value = 952.875
$ws.Cells.Item(116, 12).Value = 1006.5
$ws.Cells.Item(116, 13).Value = 1341.125
$ws.Cells.Item(116, 14).Value = -5594.5
$ws.Cells.Item(122, 8).Value = 1830.5294
$ws.Cells.Item(122, 9).Value = 1655.3077
$ws.Cells.Item(122, 11).Value = 4965.9231
$ws.Cells.Item(122, 13).Value = -2515.9231
$ws.Cells.Item(133, 8).Value = 49994.25
$ws.Cells.Item(133, 10).Value = 49992.332
$ws.Cells.Item(133, 12).Value = 49992.332
$ws.Cells.Item(133, 14).Value = -55052.332
$ws.Cells.Item(136, 8).Value = 2692.5715
$ws.Cells.Item(136, 9).Value = 1974.75
$ws.Cells.Item(136, 11).Value = 5924.25
$ws.Cells.Item(136, 13).Value = -3374.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 970.75
$ws.Cells.Item(3, 9).Value = 952.875
$ws.Cells.Item(3, 10).Value = 1006.5
$ws.Cells.Item(3, 11).Value = 952.875
$ws.Cells.Item(3, 12).Value = 1006.5
$ws.Cells.Item(3, 13).Value = -838.875
$ws.Cells.Item(3, 14).Value = -1234.5
$ws.Cells.Item(16, 8).Value = 1997.75
$ws.Cells.Item(16, 9).Value = 933.6667
$ws.Cells.Item(16, 10).Value = 5190
$ws.Cells.Item(16, 11).Value = 933.6667
$ws.Cells.Item(16, 12).Value = 5190
$ws.Cells.Item(16, 13).Value = -646.6667
$ws.Cells.Item(16, 14).Value = -5764
$ws.Cells.Item(20, 8).Value = 2355.7693
$ws.Cells.Item(20, 9).Value = 2356.3333
$ws.Cells.Item(20, 11).Value = 2356.3333
$ws.Cells.Item(20, 13).Value = -2109.3333
$ws.Cells.Item(82, 8).Value = 55662.332
$ws.Cells.Item(82, 9).Value = 6999.5
$ws.Cells.Item(82, 11).Value = 6999.5
$ws.Cells.Item(82, 13).Value = -6616.5
$ws.Cells.Item(85, 8).Value = 55662.332
$ws.Cells.Item(85, 9).Value = 6999.5
$ws.Cells.Item(85, 11).Value = 6999.5
$ws.Cells.Item(85, 13).Value = -5673.5
$ws.Cells.Item(86, 8).Value = 21388.727
$ws.Cells.Item(86, 9).Value = 16613.857
$ws.Cells.Item(86, 11).Value = 16613.857
$ws.Cells.Item(86, 13).Value = -15490.857
$ws.Cells.Item(89, 8).Value = 21388.727
$ws.Cells.Item(89, 9).Value = 16613.857
$ws.Cells.Item(89, 11).Value = 83069.285
$ws.Cells.Item(89, 13).Value = -77453.285
$ws.Cells.Item(94, 8).Value = 7017.7
$ws.Cells.Item(94, 10).Value = 14462.125
$ws.Cells.Item(94, 12).Value = 14462.125
$ws.Cells.Item(94, 14).Value = -15364.125
$ws.Cells.Item(99, 8).Value = 886
$ws.Cells.Item(99, 9).Value = 410
$ws.Cells.Item(99, 10).Value = 1124
$ws.Cells.Item(99, 11).Value = 410
$ws.Cells.Item(99, 12).Value = 1124
$ws.Cells.Item(99, 13).Value = 1088
$ws.Cells.Item(99, 14).Value = -4120
$ws.Cells.Item(105, 8).Value = 2971.818
$ws.Cells.Item(105, 9).Value = 2769.1
$ws.Cells.Item(105, 11).Value = 2769.1
$ws.Cells.Item(105, 13).Value = -1022.1
$ws.Cells.Item(107, 8).Value = 1893.625
$ws.Cells.Item(107, 9).Value = 2190.5
$ws.Cells.Item(107, 10).Value = 1003
$ws.Cells.Item(107, 11).Value = 2190.5
$ws.Cells.Item(107, 12).Value = 1003
$ws.Cells.Item(107, 13).Value = -270.5
$ws.Cells.Item(107, 14).Value = -4843
$ws.Cells.Item(124, 8).Value = 30000
$ws.Cells.Item(124, 10).Value = 30000
$ws.Cells.Item(124, 12).Value = 30000
$ws.Cells.Item(124, 14).Value = -39820
$ws.Cells.Item(134, 8).Value = 8434.979499999999
$ws.Cells.Item(134, 9).Value = 8803.076999999999
$ws.Cells.Item(134, 10).Value = 7999.9546
$ws.Cells.Item(134, 11).Value = 26409.231
$ws.Cells.Item(134, 12).Value = 23999.8638
$ws.Cells.Item(134, 13).Value = -23874.231
$ws.Cells.Item(134, 14).Value = -29069.8638
$ws.Cells.Item(140, 8).Value = 120000
$ws.Cells.Item(140, 10).Value = 120000
$ws.Cells.Item(140, 12).Value = 120000
$ws.Cells.Item(140, 14).Value = -130360
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2845.7
$ws.Cells.Item(31, 9).Value = 2174
$ws.Cells.Item(31, 10).Value = 3133.5715
$ws.Cells.Item(31, 11).Value = 2174
$ws.Cells.Item(31, 12).Value = 3133.5715
$ws.Cells.Item(31, 13).Value = -1879
$ws.Cells.Item(31, 14).Value = -3723.5715
$ws.Cells.Item(34, 8).Value = 2845.7
$ws.Cells.Item(34, 9).Value = 2174
$ws.Cells.Item(34, 10).Value = 3133.5715
$ws.Cells.Item(34, 11).Value = 2174
$ws.Cells.Item(34, 12).Value = 3133.5715
$ws.Cells.Item(34, 13).Value = -1972
$ws.Cells.Item(34, 14).Value = -3537.5715
$ws.Cells.Item(62, 8).Value = 6128.2856
$ws.Cells.Item(62, 9).Value = 3999.5
$ws.Cells.Item(62, 10).Value = 6979.8
$ws.Cells.Item(62, 11).Value = 3999.5
$ws.Cells.Item(62, 12).Value = 6979.8
$ws.Cells.Item(62, 13).Value = -3375.5
$ws.Cells.Item(62, 14).Value = -8227.799999999999
$ws.Cells.Item(65, 8).Value = 6128.2856
$ws.Cells.Item(65, 9).Value = 3999.5
$ws.Cells.Item(65, 10).Value = 6979.8
$ws.Cells.Item(65, 11).Value = 19997.5
$ws.Cells.Item(65, 12).Value = 34899
$ws.Cells.Item(65, 13).Value = -16877.5
$ws.Cells.Item(65, 14).Value = -41139
$ws.Cells.Item(94, 8).Value = 1855.2
$ws.Cells.Item(94, 9).Value = 2522
$ws.Cells.Item(94, 10).Value = 1612.7273
$ws.Cells.Item(94, 11).Value = 2522
$ws.Cells.Item(94, 12).Value = 1612.7273
$ws.Cells.Item(94, 13).Value = -2071
$ws.Cells.Item(94, 14).Value = -2514.7273
$ws.Cells.Item(99, 8).Value = 71430520
$ws.Cells.Item(99, 10).Value = 2177.5
$ws.Cells.Item(99, 12).Value = 2177.5
$ws.Cells.Item(99, 14).Value = -5173.5
$ws.Cells.Item(105, 8).Value = 1599.6
$ws.Cells.Item(105, 9).Value = 1333.3334
$ws.Cells.Item(105, 11).Value = 1333.3334
$ws.Cells.Item(105, 13).Value = 413.6666
$ws.Cells.Item(113, 8).Value = 1997.75
$ws.Cells.Item(113, 9).Value = 933.6667
$ws.Cells.Item(113, 10).Value = 5190
$ws.Cells.Item(113, 11).Value = 933.6667
$ws.Cells.Item(113, 12).Value = 5190
$ws.Cells.Item(113, 13).Value = 1236.3333
$ws.Cells.Item(113, 14).Value = -9530
$ws.Cells.Item(126, 8).Value = 71430520
$ws.Cells.Item(126, 10).Value = 2177.5
$ws.Cells.Item(126, 12).Value = 6532.5
$ws.Cells.Item(126, 14).Value = -11472.5
$ws.Cells.Item(134, 8).Value = 3281.0715
$ws.Cells.Item(134, 9).Value = 3326.75
$ws.Cells.Item(134, 11).Value = 9980.25
$ws.Cells.Item(134, 13).Value = -7445.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 4015.5
$ws.Cells.Item(3, 9).Value = 4015.5
$ws.Cells.Item(3, 11).Value = 12046.5
$ws.Cells.Item(3, 13).Value = -11934.5
$ws.Cells.Item(17, 8).Value = 367.625
$ws.Cells.Item(17, 9).Value = 115.25
$ws.Cells.Item(17, 10).Value = 620
$ws.Cells.Item(17, 11).Value = 345.75
$ws.Cells.Item(17, 12).Value = 1860
$ws.Cells.Item(17, 13).Value = -176.75
$ws.Cells.Item(17, 14).Value = -2198
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).Value = 0
$ws.Cells.Item(81, 8).Value = 1855.4286
$ws.Cells.Item(81, 9).Value = 1855.4286
$ws.Cells.Item(81, 11).Value = 5566.2858
$ws.Cells.Item(81, 13).Value = -4443.2858
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).Value = 0
$ws.Cells.Item(84, 8).Value = 1855.4286
$ws.Cells.Item(84, 9).Value = 1855.4286
$ws.Cells.Item(84, 11).Value = 16698.8574
$ws.Cells.Item(84, 13).Value = -11082.8574
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(92, 14).Value = 0
$ws.Cells.Item(114, 8).Value = 1268.5714
$ws.Cells.Item(114, 10).Value = 1268.5714
$ws.Cells.Item(114, 12).Value = 3805.7142
$ws.Cells.Item(114, 14).Value = -10313.7142
$ws.Cells.Item(117, 8).Value = 3593.9
$ws.Cells.Item(117, 9).Value = 2633.3333
$ws.Cells.Item(117, 10).Value = 4005.5715
$ws.Cells.Item(117, 11).Value = 7899.999899999999
$ws.Cells.Item(117, 12).Value = 12016.7145
$ws.Cells.Item(117, 13).Value = -4457.999899999999
$ws.Cells.Item(117, 14).Value = -18900.7145
$ws.Cells.Item(125, 8).Value = 22142.857
$ws.Cells.Item(132, 8).Value = 3298.353
$ws.Cells.Item(132, 9).Value = 1998.5
$ws.Cells.Item(132, 10).Value = 3698.3076
$ws.Cells.Item(132, 11).Value = 17986.5
$ws.Cells.Item(132, 12).Value = 33284.7684
$ws.Cells.Item(132, 13).Value = -15456.5
$ws.Cells.Item(132, 14).Value = -38344.7684

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 258675.88
$ws.Cells.Item(22, 10).Value = 15000
$ws.Cells.Item(22, 12).Value = 15000
$ws.Cells.Item(22, 14).Value = -16058
$ws.Cells.Item(97, 8).Value = 2181.4285
$ws.Cells.Item(97, 9).Value = 1295.875
$ws.Cells.Item(97, 10).Value = 2726.3845
$ws.Cells.Item(97, 11).Value = 1295.875
$ws.Cells.Item(97, 12).Value = 2726.3845
$ws.Cells.Item(97, 13).Value = -799.875
$ws.Cells.Item(97, 14).Value = -3718.3845
$ws.Cells.Item(102, 8).Value = 2890
$ws.Cells.Item(102, 9).Value = 2605.2
$ws.Cells.Item(102, 10).Value = 3500.2856
$ws.Cells.Item(102, 11).Value = 2605.2
$ws.Cells.Item(102, 12).Value = 3500.2856
$ws.Cells.Item(102, 13).Value = -983.1999999999998
$ws.Cells.Item(102, 14).Value = -6744.2856
$ws.Cells.Item(103, 8).Value = 61249.5
$ws.Cells.Item(103, 10).Value = 61249.5
$ws.Cells.Item(103, 12).Value = 61249.5
$ws.Cells.Item(103, 14).Value = -63593.5
$ws.Cells.Item(122, 8).Value = 2800.125
$ws.Cells.Item(122, 9).Value = 2692.3333
$ws.Cells.Item(122, 10).Value = 3123.5
$ws.Cells.Item(122, 11).Value = 8076.999899999999
$ws.Cells.Item(122, 12).Value = 9370.5
$ws.Cells.Item(122, 13).Value = -5626.999899999999
$ws.Cells.Item(122, 14).Value = -14270.5
$ws.Cells.Item(126, 8).Value = 5986.25
$ws.Cells.Item(126, 9).Value = 2965.3333
$ws.Cells.Item(126, 10).Value = 7798.8
$ws.Cells.Item(126, 11).Value = 8895.999899999999
$ws.Cells.Item(126, 12).Value = 23396.4
$ws.Cells.Item(126, 13).Value = -6425.999899999999
$ws.Cells.Item(126, 14).Value = -28336.4
$ws.Cells.Item(132, 8).Value = 2176.718
$ws.Cells.Item(132, 9).Value = 1691.2
$ws.Cells.Item(132, 10).Value = 6425
$ws.Cells.Item(132, 11).Value = 5073.6
$ws.Cells.Item(132, 12).Value = 19275
$ws.Cells.Item(132, 13).Value = -2543.6
$ws.Cells.Item(132, 14).Value = -24335

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2995.8333
$ws.Cells.Item(7, 9).Value = 2994
$ws.Cells.Item(7, 11).Value = 2994
$ws.Cells.Item(7, 13).Value = -2882
$ws.Cells.Item(16, 8).Value = 2665.3333
$ws.Cells.Item(16, 9).Value = 2665.3333
$ws.Cells.Item(16, 11).Value = 2665.3333
$ws.Cells.Item(16, 13).Value = -2495.3333
$ws.Cells.Item(22, 8).Value = 959.8
$ws.Cells.Item(22, 9).Value = 500
$ws.Cells.Item(22, 11).Value = 500
$ws.Cells.Item(22, 13).Value = -205
$ws.Cells.Item(27, 8).Value = 959.8
$ws.Cells.Item(27, 9).Value = 500
$ws.Cells.Item(27, 11).Value = 500
$ws.Cells.Item(27, 13).Value = -393
$ws.Cells.Item(40, 8).Value = 2314.3333
$ws.Cells.Item(40, 9).Value = 2204.875
$ws.Cells.Item(40, 11).Value = 2204.875
$ws.Cells.Item(40, 13).Value = -2068.875
$ws.Cells.Item(46, 8).Value = 2042.7142
$ws.Cells.Item(46, 10).Value = 2674.75
$ws.Cells.Item(46, 12).Value = 2674.75
$ws.Cells.Item(46, 14).Value = -3050.75
$ws.Cells.Item(82, 8).Value = 2285.75
$ws.Cells.Item(82, 9).Value = 2100.625
$ws.Cells.Item(82, 11).Value = 2100.625
$ws.Cells.Item(82, 13).Value = -1739.625
$ws.Cells.Item(85, 8).Value = 2285.75
$ws.Cells.Item(85, 9).Value = 2100.625
$ws.Cells.Item(85, 11).Value = 2100.625
$ws.Cells.Item(85, 13).Value = -852.625
$ws.Cells.Item(126, 8).Value = 2995.8333
$ws.Cells.Item(126, 9).Value = 2994
$ws.Cells.Item(126, 11).Value = 8982
$ws.Cells.Item(126, 13).Value = -6512
$ws.Cells.Item(136, 8).Value = 2510.125
$ws.Cells.Item(136, 9).Value = 2205.3
$ws.Cells.Item(136, 10).Value = 3018.1667
$ws.Cells.Item(136, 11).Value = 6615.900000000001
$ws.Cells.Item(136, 12).Value = 9054.500100000001
$ws.Cells.Item(136, 13).Value = -4065.900000000001
$ws.Cells.Item(136, 14).Value = -14154.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 34349
$ws.Cells.Item(51, 9).Value = 27970
$ws.Cells.Item(51, 10).Value = 37538.5
$ws.Cells.Item(51, 11).Value = 27970
$ws.Cells.Item(51, 12).Value = 37538.5
$ws.Cells.Item(51, 13).Value = -27460
$ws.Cells.Item(51, 14).Value = -38558.5
$ws.Cells.Item(54, 8).Value = 50000
$ws.Cells.Item(54, 10).Value = 60000
$ws.Cells.Item(54, 12).Value = 60000
$ws.Cells.Item(54, 14).Value = -61040
$ws.Cells.Item(100, 8).Value = 4166
$ws.Cells.Item(100, 9).Value = 3900
$ws.Cells.Item(100, 10).Value = 4299
$ws.Cells.Item(100, 11).Value = 7800
$ws.Cells.Item(100, 12).Value = 8598
$ws.Cells.Item(100, 13).Value = -7259
$ws.Cells.Item(100, 14).Value = -9680
$ws.Cells.Item(107, 8).Value = 923.4545000000001
$ws.Cells.Item(107, 9).Value = 889.125
$ws.Cells.Item(107, 11).Value = 2667.375
$ws.Cells.Item(107, 13).Value = -747.375
$ws.Cells.Item(122, 8).Value = 5459.9473
$ws.Cells.Item(122, 9).Value = 5814.0586
$ws.Cells.Item(122, 11).Value = 17442.1758
$ws.Cells.Item(122, 13).Value = -14992.1758
$ws.Cells.Item(126, 8).Value = 3459
$ws.Cells.Item(126, 9).Value = 1534.3
$ws.Cells.Item(126, 10).Value = 6208.5713
$ws.Cells.Item(126, 11).Value = 4602.9
$ws.Cells.Item(126, 12).Value = 18625.7139
$ws.Cells.Item(126, 13).Value = -2132.9
$ws.Cells.Item(126, 14).Value = -23565.7139
$ws.Cells.Item(136, 8).Value = 6198.914
$ws.Cells.Item(136, 9).Value = 5302.2144
$ws.Cells.Item(136, 10).Value = 9785.714
$ws.Cells.Item(136, 11).Value = 15906.6432
$ws.Cells.Item(136, 12).Value = 29357.142
$ws.Cells.Item(136, 13).Value = -13356.6432
$ws.Cells.Item(136, 14).Value = -34457.142
